$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/20/2025  Through  10/26/2025"

# --- Cells that change value but keep their original type/style (number <-> number) ---
$ws.Range("M14").Value = -66.666666666666
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = 5.882352941176
$ws.Range("L15").Value = 63.636363636363
$ws.Range("N15").Value = 12.5
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 81.818181818181
$ws.Range("I16").Value = 160
$ws.Range("J16").Value = 129
$ws.Range("K16").Value = 24.031007751938
$ws.Range("L16").Value = 7.382550335570
$ws.Range("M16").Value = 9.589041095890
$ws.Range("N16").Value = -81.330221703617
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -11.764705882352
$ws.Range("I17").Value = 242
$ws.Range("J17").Value = 195
$ws.Range("K17").Value = 24.102564102564
$ws.Range("L17").Value = 26.041666666666
$ws.Range("M17").Value = 108.620689655172
$ws.Range("N17").Value = 2.542372881355
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 11.764705882352
$ws.Range("I18").Value = 303
$ws.Range("J18").Value = 164
$ws.Range("K18").Value = 84.756097560975
$ws.Range("L18").Value = 63.783783783783
$ws.Range("M18").Value = 46.376811594202
$ws.Range("N18").Value = -77.235161532682
$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = 37.5
$ws.Range("G19").Value = 78
$ws.Range("H19").Value = 56.410256410256
$ws.Range("I19").Value = 1004
$ws.Range("J19").Value = 857
$ws.Range("K19").Value = 17.152858809801
$ws.Range("L19").Value = 16.744186046511
$ws.Range("M19").Value = -15.559293523969
$ws.Range("N19").Value = -55.298308103294
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -25
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = 19.444444444444
$ws.Range("L20").Value = -20.370370370370
$ws.Range("M20").Value = 7.5
$ws.Range("N20").Value = -95.784313725490
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = 7.317073170731
$ws.Range("F21").Value = 182
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = 40
$ws.Range("I21").Value = 1771
$ws.Range("J21").Value = 1400
$ws.Range("K21").Value = 26.5
$ws.Range("L21").Value = 21.969696969697
$ws.Range("M21").Value = 3.567251461988
$ws.Range("N21").Value = -68.995098039215
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -66.666666666666
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 63
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = 21.153846153846
$ws.Range("L22").Value = -18.181818181818
$ws.Range("M22").Value = 12.5
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = 37.5
$ws.Range("L23").Value = 10
$ws.Range("M23").Value = -26.666666666666
$ws.Range("C24").Value = 101
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = 77.192982456140
$ws.Range("F24").Value = 279
$ws.Range("G24").Value = 259
$ws.Range("H24").Value = 7.722007722007
$ws.Range("I24").Value = 2293
$ws.Range("J24").Value = 2439
$ws.Range("K24").Value = -5.986059860598
$ws.Range("L24").Value = 27.036011080332
$ws.Range("M24").Value = 50.558108995403
$ws.Range("C25").Value = 92
$ws.Range("D25").Value = 47
$ws.Range("E25").Value = 95.744680851063
$ws.Range("F25").Value = 249
$ws.Range("G25").Value = 208
$ws.Range("H25").Value = 19.711538461538
$ws.Range("I25").Value = 1877
$ws.Range("J25").Value = 2070
$ws.Range("K25").Value = -9.323671497584
$ws.Range("L25").Value = 39.865871833084
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = -6.521739130434
$ws.Range("I26").Value = 500
$ws.Range("J26").Value = 497
$ws.Range("K26").Value = 0.603621730382
$ws.Range("L26").Value = 19.617224880382
$ws.Range("M26").Value = 42.045454545454
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = -24.137931034482
$ws.Range("L27").Value = 4.761904761904
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 105
$ws.Range("J28").Value = 109
$ws.Range("K28").Value = -3.669724770642
$ws.Range("L28").Value = 7.142857142857
$ws.Range("M29").Value = -66.666666666666
$ws.Range("M30").Value = -50
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = -27.272727272727

# --- Cells that change between shared-string placeholder ("0"/"***.*") and numeric value ---
# These require both a value change and a style/type change. We set the raw value first
# (using a leading apostrophe to force text where needed), then copy formatting only from a
# donor cell with the correct target style so the cell xf index matches what Excel would
# naturally use (s="13" for text placeholders, s="14" for plain integers, s="15" for percentages).

$ws.Range("C15").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C20").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("G23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G23").PasteSpecial(-4122)

$ws.Range("H23").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H23").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D31").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("E31").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$ws.Range("G31").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("H31").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$excel.CutCopyMode = 0
